$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells G1, H1 - copy style from existing header (F1) then set text
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Updated B/C/D values for rows 2-10
$ws.Range("B2").Value = 1.294123935468054
$ws.Range("C2").Value = 0.7891738527757535
$ws.Range("D2").Value = 0.9563752695328483

$ws.Range("B3").Value = 6.926509703753763
$ws.Range("C3").Value = 0.9010492682785619
$ws.Range("D3").Value = 1.957121876262595

$ws.Range("B4").Value = 4.245823787260473
$ws.Range("C4").Value = 0.7902735092942522
$ws.Range("D4").Value = 1.627343786371562

$ws.Range("B5").Value = 4.041165657826554
$ws.Range("C5").Value = 0.9973455951643054
$ws.Range("D5").Value = 1.695698936855346

$ws.Range("B6").Value = 2.867831406231114
$ws.Range("C6").Value = 0.9687500485729165
$ws.Range("D6").Value = 1.436044889417979

$ws.Range("B7").Value = 2.470467071850259
$ws.Range("C7").Value = 0.9986148863130003
$ws.Range("D7").Value = 1.277382901525897

$ws.Range("B8").Value = 2.016866106204859
$ws.Range("C8").Value = 0.9976897741085226
$ws.Range("D8").Value = 1.14505494204693

$ws.Range("B9").Value = 15.77119385868702
$ws.Range("C9").Value = 0.8114662754641242
$ws.Range("D9").Value = 3.132059411705479

$ws.Range("B10").Value = 2.029428910754003
$ws.Range("C10").Value = 0.9939503352753355
$ws.Range("D10").Value = 1.134732011086487

# New G/H columns for rows 2-10 (Elapsed Time / CPU)
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 7).Value = 1.127317944850074
    $ws.Cells.Item($r, 8).Value = 0.985
}
